{"js": "// Drop the \".io\" suffix from the \"Zeromon.io\" line (leaving just \"Zeromon\"),\n// and change \"Amazon Alexa Skill for SeizureTracker(.com)\" to\n// \"Amazon Alexa Skill for SeizureTracker.com\".\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nfor (const p of paragraphs.items) {\n  if (p.text === \"Zeromon.io\") {\n    const hits = p.search(\".io\", { matchCase: true });\n    hits.load(\"items\");\n    await context.sync();\n    if (hits.items.length > 0) {\n      hits.items[hits.items.length - 1].delete();\n    }\n  } else if (p.text === \"Amazon Alexa Skill for SeizureTracker(.com)\") {\n    const hits = p.search(\"(.com)\", { matchCase: true });\n    hits.load(\"items\");\n    await context.sync();\n    if (hits.items.length > 0) {\n      hits.items[0].insertText(\".com\", \"Replace\");\n    }\n  }\n}\n\nawait context.sync();\n", "ps1": "# Drop the \".io\" suffix from the \"Zeromon.io\" line (leaving just \"Zeromon\"),\n# and change \"Amazon Alexa Skill for SeizureTracker(.com)\" to\n# \"Amazon Alexa Skill for SeizureTracker.com\".\n\n$d = $word.ActiveDocument\n\n# wdFindContinue = 1, wdReplaceOne = 2\n$wdFindContinue = 1\n$wdReplaceOne = 2\n\n$r1 = $d.Content\n$r1.Find.Execute(\"Zeromon.io\", $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, \"Zeromon\", $wdReplaceOne) | Out-Null\n\n$r2 = $d.Content\n$r2.Find.Execute(\"SeizureTracker(.com)\", $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, \"SeizureTracker.com\", $wdReplaceOne) | Out-Null\n"}
